$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trim the trailing space off the "facoep " value (database param, Usar=True row)
$ws.Range("B6").Value = "facoep"

# Reflect the user's click/selection on B7 after the edit (also nudges the
# visible scroll position so row 7 is in view, matching topLeftCell="B1")
$ws.Application.Goto($ws.Range("B7"))
$ws.Range("B7").Select()
try {
    $win = $excel.ActiveWindow
    $win.ScrollColumn = 2
    $win.ScrollRow = 1
} catch {
    # Scroll-position persistence isn't guaranteed by every host; ignore.
}
